# Apply the cryptos list update (cell-value changes only; no style changes intended).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values: assign directly.
$plainUpdates = @{
    "D2" = "96.522.85"
    "E2" = "  -1.39%  "
    "D3" = "3.331.90"
    "E3" = "  -2.79%  "
    "E4" = "  -0.04%  "
    "E5" = "  -2.42%  "
    "E6" = "  -0.49%  "
    "E7" = "  -6.72%  "
    "E8" = "  -2.55%  "
    "E10" = "  -6.39%  "
    "D11" = "3.331.58"
    "E11" = "  -2.70%  "
    "E12" = "  -3.60%  "
    "E13" = "  -4.21%  "
    "D14" = "96.223.73"
    "E14" = "  -1.51%  "
    "E15" = "  -4.17%  "
    "E16" = "  -3.80%  "
    "D17" = "3.956.12"
    "E17" = "  -2.61%  "
    "E18" = "  -3.46%  "
    "D19" = "3.329.83"
    "E19" = "  -2.73%  "
    "E20" = "  -2.70%  "
    "E21" = "  +1.84%  "
    "B22" = "BitcoinCash"
    "C22" = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
    "E22" = "  -1.06%  "
    "B23" = "Uniswap"
    "C23" = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
    "E23" = "  -4.78%  "
    "E24" = "  -3.32%  "
    "E26" = "  +6.35%  "
    "E27" = "  -3.06%  "
    "E28" = "  -6.12%  "
    "D29" = "3.507.00"
    "E29" = "  -2.55%  "
    "E30" = "  -8.78%  "
    "E31" = "  +0.01%  "
    "E32" = "  -4.10%  "
    "E33" = "  -6.39%  "
    "E34" = "  +9.26%  "
    "E35" = "  -0.08%  "
    "E36" = "  -4.93%  "
    "E37" = "  -7.13%  "
    "E38" = "  +3.62%  "
    "E39" = "  -3.31%  "
    "E40" = "  +0.00%  "
    "E41" = "  -3.79%  "
    "E42" = "  -2.13%  "
    "E43" = "  -1.64%  "
    "E44" = "  -3.99%  "
    "E45" = "  +1.01%  "
    "E46" = "  -1.28%  "
    "E47" = "  +5.55%  "
    "E48" = "  +0.41%  "
    "E49" = "  +1.28%  "
    "E50" = "  +3.59%  "
    "E51" = "  -5.53%  "
}
foreach ($ref in $plainUpdates.Keys) {
    $ws.Range($ref).Value = $plainUpdates[$ref]
}

# Values that look like plain numbers (e.g. "250.04") must stay as literal text
# (matching the source price-string formatting), so force Text format before
# assigning, then restore the default style so no stray formatting is left behind.
$numericLookingUpdates = @{
    "D5" = "250.04"
    "D6" = "655.14"
    "D13" = "40.27"
    "D15" = "6.09"
    "D18" = "8.52"
    "D20" = "17.09"
    "D21" = "0.523"
    "D22" = "504.76"
    "D23" = "10.54"
    "D26" = "6.58"
    "D27" = "96.41"
    "D28" = "12.08"
    "D32" = "11.06"
    "D37" = "27.98"
    "D42" = "508.89"
    "D44" = "0.832"
    "D45" = "0.0422"
    "D47" = "1.67"
    "D48" = "5.50"
    "D49" = "8.36"
    "D50" = "53.27"
}
foreach ($ref in $numericLookingUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $numericLookingUpdates[$ref]
    $cell.Style = "Normal"
}

Write-Output "Applied $($plainUpdates.Count + $numericLookingUpdates.Count) cell updates"
